# "updated alex scheme and some gui refinements"
#
# Data update: two IP addresses in the alex_scheme sheet changed:
#   - Row 49 (V 2299 Ant. 2):            169.254.1.53 -> 169.254.1.59
#   - Row 51 (Mobinil Nubariya Rd Ant. 1): 169.254.1.68 -> 169.254.1.60
# Columns B, C and D hold the raw IP text, column G repeats it (E/F hold
# the antenna name instead, so they are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B49").Value = "169.254.1.59"
$ws.Range("C49").Value = "169.254.1.59"
$ws.Range("D49").Value = "169.254.1.59"
$ws.Range("G49").Value = "169.254.1.59"

$ws.Range("B51").Value = "169.254.1.60"
$ws.Range("C51").Value = "169.254.1.60"
$ws.Range("D51").Value = "169.254.1.60"
$ws.Range("G51").Value = "169.254.1.60"

# GUI refinement: move the selection/viewport down to where the edits were
# made instead of leaving it parked at D8.
$ws.Range("G49").Select()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
